$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7741817235946655
$ws.Range("B1").Value = 1.44614851474762
$ws.Range("D1").Value = 1.85819149017334
$ws.Range("E1").Value = 1.188699126243591
